$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A56").Value = "you are not a man yet"
$ws.Range("A60").Value = "Can you send me the email from last week"

[void]$ws.Range("H52").Select()
